# Generate Report for Archive
#
# Update the localization status text from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview!E2:F3 and the "Status"
# column on the per-locale sheets), then shrink the now-narrower Status
# columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overviewUsed = $overview.UsedRange
for ($r = 1; $r -le $overviewUsed.Rows.Count; $r++) {
    foreach ($c in 5, 6) {
        $cell = $overview.Cells.Item($r, $c)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- zh-cn / de-de sheets: column C holds the status text ---
foreach ($ws in $zhcn, $dede) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# Resize the Status columns to fit the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
